$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.891.60'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -5.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.274.37'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -6.02%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '558.78'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.53'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.267.12'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -9.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.584'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.29'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -8.05%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -7.22%  '
$ws.Range('B14').Value = 'BitcoinCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '635.02'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.61'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.799.62'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.857.72'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.23%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.274.79'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.91%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -8.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.904'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.19'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '107.39'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +8.66%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -7.03%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -7.52%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -7.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.51'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.68'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.27'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.07%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.49%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -7.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.99'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.26%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.56'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.699.04'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '522.43'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.37'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0728'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.62%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.69'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '32.76'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.32'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.07%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -9.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.28'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.55%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -6.48%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.23%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -8.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.78%  '
